$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data corrections, rows 176-181 (existing rows, only the input
# columns change; the B/H/J/K formula columns recompute automatically) ---

# Row 176: nouveaux cas positifs 3 -> 4
$ws.Range("C176").Value = 4

# Row 177: nouveaux cas positifs 8 -> 7
$ws.Range("C177").Value = 7

# Row 178: hospitalisations hors SI 8 -> 7 ; nouvelles sorties 0 -> 1
$ws.Range("G178").Value = 7
$ws.Range("I178").Value = 1

# Row 179: nouveaux cas positifs 10 -> 11 ; hospitalisations hors SI 8 -> 7
$ws.Range("C179").Value = 11
$ws.Range("G179").Value = 7

# Row 180: nouveaux cas positifs 2 -> 3 ; hospitalisations hors SI 8 -> 7
$ws.Range("C180").Value = 3
$ws.Range("G180").Value = 7

# Row 181: nouveaux cas positifs 0 -> 7 ; hospitalisations hors SI 8 -> 6 ;
# nouvelles sorties 0 -> 1
$ws.Range("C181").Value = 7
$ws.Range("G181").Value = 6
$ws.Range("I181").Value = 1

# --- Row 182: was a blank placeholder row (formulas present but no
# inputs, so they evaluated to ""); fill in the day's real figures ---
$ws.Range("C182").Value = 1
$ws.Range("D182").Value = 0
$ws.Range("E182").Value = 1
$ws.Range("F182").Value = 1
$ws.Range("G182").Value = 6
$ws.Range("I182").Value = 0
$ws.Range("L182").Value = "0"
$ws.Range("M182").Value = "0"

# --- View state: move the active cell down to where the new data was
# entered (frozen header rows 1-2 stay put) ---
$ws.Range("I184").Select()
